# Applies the recorded edit: within each of three row-groups, the
# "observation" data (taxon id, names, coordinates, etc.) is rotated one
# row down, with the last row's data wrapping around to the first row.
# Row-level metadata that is identical across a whole group (locality,
# county, date, reporter, ...) is left untouched.
#
#   group {29,30,31}: new29 = old31, new30 = old29, new31 = old30
#   group {34,35,36}: new34 = old36, new35 = old34, new36 = old35
#   group {52,53}   : new52 = old53, new53 = old52

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that actually carry per-observation data in these rows.
$cols = @("A","B","D","E","F","G","H","Q","R","AC","AJ","AK","AO")

function Get-RowSnapshot($row) {
    $snap = @{}
    foreach ($col in $cols) {
        $snap[$col] = $ws.Range("$col$row").Value2
    }
    return $snap
}

function Set-RowFromSnapshot($row, $snap) {
    foreach ($col in $cols) {
        $val = $snap[$col]
        if ($val -eq $null) {
            $ws.Range("$col$row").ClearContents()
        } else {
            $ws.Range("$col$row").Value = $val
        }
    }
}

# --- capture every source row first (reads must happen before any writes) ---
$snap29 = Get-RowSnapshot 29
$snap30 = Get-RowSnapshot 30
$snap31 = Get-RowSnapshot 31

$snap34 = Get-RowSnapshot 34
$snap35 = Get-RowSnapshot 35
$snap36 = Get-RowSnapshot 36

$snap52 = Get-RowSnapshot 52
$snap53 = Get-RowSnapshot 53

# --- write the rotated data back ---
Set-RowFromSnapshot 29 $snap31
Set-RowFromSnapshot 30 $snap29
Set-RowFromSnapshot 31 $snap30

Set-RowFromSnapshot 34 $snap36
Set-RowFromSnapshot 35 $snap34
Set-RowFromSnapshot 36 $snap35

Set-RowFromSnapshot 52 $snap53
Set-RowFromSnapshot 53 $snap52
